# Update "carjacking arrests by month, year-over-year" workbook to add data
# for 2021-12-18 (commit message says "2021-12-26" but the diff itself shows
# the sheet/date label moving from 12-17 to 12-18, so we follow the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the (only) worksheet/tab to reflect the new "through" date.
$ws.Name = "Through 2021-12-18"

# Update the label cell for the December row.
$ws.Range("A14").Value = "December (through 12-18)"

# --- Row 14: "December (through 12-18)" ---
# 2015
$ws.Range("C14").Value = 20
$ws.Range("D14").Value = 0.1304
# 2016
$ws.Range("F14").Value = 52
$ws.Range("G14").Value = 0.0877
# 2017
$ws.Range("H14").Value = 8
$ws.Range("I14").Value = 62
$ws.Range("J14").Value = 0.1143
# 2018
$ws.Range("K14").Value = 4
$ws.Range("L14").Value = 38
$ws.Range("M14").Value = 0.0952
# 2019
$ws.Range("O14").Value = 26
$ws.Range("P14").Value = 0.1034
# 2020
$ws.Range("R14").Value = 79
$ws.Range("S14").Value = 0.0482
# 2021
$ws.Range("U14").Value = 128
$ws.Range("V14").Value = 0.0154

# --- Row 15: "Total" ---
# 2015
$ws.Range("C15").Value = 278
$ws.Range("D15").Value = 0.1146
# 2016
$ws.Range("F15").Value = 556
$ws.Range("G15").Value = 0.1032
# 2017
$ws.Range("H15").Value = 71
$ws.Range("I15").Value = 820
$ws.Range("J15").Value = 0.0797
# 2018
$ws.Range("K15").Value = 78
$ws.Range("L15").Value = 646
$ws.Range("M15").Value = 0.1077
# 2019
$ws.Range("O15").Value = 506
$ws.Range("P15").Value = 0.1012
# 2020
$ws.Range("R15").Value = 1279
$ws.Range("S15").Value = 0.0505
# 2021
$ws.Range("U15").Value = 1670
$ws.Range("V15").Value = 0.0581
